# Generate Report for Handback
# This script updates the handback-status workbook so that the two
# e2e/*.md rows now refer to the newly generated guids/hashes and
# timestamps, mirroring a fresh handback-report generation.

$wb = $excel.ActiveWorkbook

$oldUuid1 = "3dcebe8c-00c0-4864-8ebe-cfad31e4405c"
$oldUuid2 = "869b18e7-eb08-4135-b7d6-2f61e966e224"
$newUuid1 = "43ba885b-911d-48b2-847f-0c9ff3369a15"
$newUuid2 = "ffff6cfd71fe-2663-4174-9774-60d3076beb11"
$newHash  = "1788a6d4d7b073fbd3126a615837aab506dbce54"

# -------------------------------------------------------------------
# Sheet "Overview"
# -------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newUuid1.md"
$wsOverview.Range("B2").Value = "e2e\$newUuid1.md"
$wsOverview.Range("G2").Value = "2016-09-03 17:10:37"

$wsOverview.Range("A3").Value = "$newUuid2.md"
$wsOverview.Range("B3").Value = "e2e\$newUuid2.md"
$wsOverview.Range("G3").Value = "2016-09-03 17:10:37"

# Recreate the hyperlinks on B2/B3 so their display text matches the
# renamed file names while the link targets (addresses) stay the same.
$overviewLinkAddr2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bfce9f0d6bd0d8a2bc145b121a9baab849c3f089/e2e/$oldUuid1.md"
$overviewLinkAddr3 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bfce9f0d6bd0d8a2bc145b121a9baab849c3f089/e2e/$oldUuid2.md"

$wsOverview.Range("A1").Hyperlinks.Delete()

$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $overviewLinkAddr2, [Type]::Missing, [Type]::Missing, "e2e\$newUuid1.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $overviewLinkAddr3, [Type]::Missing, [Type]::Missing, "e2e\$newUuid2.md") | Out-Null

# -------------------------------------------------------------------
# Sheet "zh-cn"
# -------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newUuid1.md"
$wsZhCn.Range("G2").Value = "$newUuid1.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-09-03 17:10:32"
$wsZhCn.Range("I2").Value = "$newUuid1.md"
$wsZhCn.Range("J2").Value = "$newUuid1.$newHash.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-09-03 17:10:49"

$wsZhCn.Range("A3").Value = "$newUuid2.md"
$wsZhCn.Range("G3").Value = "$newUuid1.$newHash.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-03 17:10:32"
$wsZhCn.Range("I3").Value = "$newUuid2.md"
$wsZhCn.Range("J3").Value = "$newUuid1.$newHash.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-09-03 17:10:49"

$zhcnLinkAddrA2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bfce9f0d6bd0d8a2bc145b121a9baab849c3f089/e2e/$oldUuid1.md"
$zhcnLinkAddrI2 = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/db14f8af43412e7c15348464fcb8c73b5eccf3b6/e2e/$oldUuid1.md"
$zhcnLinkAddrA3 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bfce9f0d6bd0d8a2bc145b121a9baab849c3f089/e2e/$oldUuid2.md"
$zhcnLinkAddrI3 = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/db14f8af43412e7c15348464fcb8c73b5eccf3b6/e2e/$oldUuid2.md"

$wsZhCn.Range("A1").Hyperlinks.Delete()

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $zhcnLinkAddrA2, [Type]::Missing, [Type]::Missing, "$newUuid1.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $zhcnLinkAddrI2, [Type]::Missing, [Type]::Missing, "$newUuid1.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $zhcnLinkAddrA3, [Type]::Missing, [Type]::Missing, "$newUuid2.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $zhcnLinkAddrI3, [Type]::Missing, [Type]::Missing, "$newUuid2.md") | Out-Null

# -------------------------------------------------------------------
# Sheet "de-de"
# -------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newUuid1.md"
$wsDeDe.Range("G2").Value = "$newUuid1.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-09-03 17:10:37"
$wsDeDe.Range("I2").Value = "$newUuid1.md"
$wsDeDe.Range("J2").Value = "$newUuid1.$newHash.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-03 17:10:57"

$wsDeDe.Range("A3").Value = "$newUuid2.md"
$wsDeDe.Range("G3").Value = "$newUuid1.$newHash.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-03 17:10:37"
$wsDeDe.Range("I3").Value = "$newUuid2.md"
$wsDeDe.Range("J3").Value = "$newUuid1.$newHash.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-09-03 17:10:57"

$dedeLinkAddrA2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bfce9f0d6bd0d8a2bc145b121a9baab849c3f089/e2e/$oldUuid1.md"
$dedeLinkAddrI2 = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/7bfc6021ab40e21968a892393d0c821934fad151/e2e/$oldUuid1.md"
$dedeLinkAddrA3 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bfce9f0d6bd0d8a2bc145b121a9baab849c3f089/e2e/$oldUuid2.md"
$dedeLinkAddrI3 = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/7bfc6021ab40e21968a892393d0c821934fad151/e2e/$oldUuid2.md"

$wsDeDe.Range("A1").Hyperlinks.Delete()

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $dedeLinkAddrA2, [Type]::Missing, [Type]::Missing, "$newUuid1.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $dedeLinkAddrI2, [Type]::Missing, [Type]::Missing, "$newUuid1.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $dedeLinkAddrA3, [Type]::Missing, [Type]::Missing, "$newUuid2.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $dedeLinkAddrI3, [Type]::Missing, [Type]::Missing, "$newUuid2.md") | Out-Null
